$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Event rows data: event_id, fecha, jugador_A, jugador_B, pronostico, cuota
$rows = @(
    @(14633444, "2025-09-08", "Dennis Novak", "Sandro Kopp", "Gana Sandro Kopp", 2.1),
    @(14633443, "2025-09-08", "Eric Vanshelboim", "Jelle Sels", "Gana Jelle Sels", 1.83),
    @(14633183, "2025-09-08", "Kilian Feldbausch", "Alexander Weis", "Gana Alexander Weis", 3.25),
    @(14631272, "2025-09-08", "Yuta Shimizu", "Juan Manuel Cerundolo", "Gana Juan Manuel Cerundolo", 1.67),
    @(14637089, "2025-09-08", "Marvin Möller", "Jozef Kovalik", "Gana Jozef Kovalik", 1.83),
    @(14636896, "2025-09-08", "Tadeas Paroulek", "Lorenzo Bocchi", "Gana Lorenzo Bocchi", 3),
    @(14637455, "2025-09-08", "Alexey Vatutin", "Niels Visker", "Gana Niels Visker", 4.5),
    @(14637123, "2025-09-08", "Rudolf Molleker", "Maik Steiner", "Gana Maik Steiner", 4.5),
    @(14637454, "2025-09-08", "Jakub Nicod", "Michael Vrbensky", "Gana Michael Vrbensky", 2.62),
    @(14637031, "2025-09-08", "Enzo Couacaud", "Robin Catry", "Gana Robin Catry", 2.62),
    @(14637438, "2025-09-08", "Maxime Janvier", "Adrien Gobat", "Gana Adrien Gobat", 2.1),
    @(14637071, "2025-09-08", "Pol Martin Tiffon", "Michele Ribecai", "Gana Michele Ribecai", 2.5),
    @(14631300, "2025-09-08", "Kasidit Samrej", "Linghao Zhang", "Gana Linghao Zhang", 8),
    @(14637056, "2025-09-08", "Luca Castelnuovo", "Ye Cong Mo", "Gana Ye Cong Mo", 2.38),
    @(14637057, "2025-09-08", "Alexandr Binda", "Evan Zhu", "Gana Evan Zhu", 2.62)
)

$startRow = 2
$endRow = $startRow + $rows.Count - 1

# Force column B (fecha) to be treated as text so date strings are not
# auto-converted into date serials by the COM layer.
$ws.Range("B$startRow`:B$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

# Remove the temporary text format from column B now that the values
$ws.Range("B$startRow`:B$endRow").ClearFormats()

# Columns G (resultado) and H (profit) are present but empty for each new row.
$ws.Range("G$startRow`:H$endRow").NumberFormat = "@"
$ws.Range("G$startRow`:H$endRow").Value = ""
$ws.Range("G$startRow`:H$endRow").ClearFormats()

